# Add "(분류모델)" after the existing "예측" run in the title placeholder
# of slide 1, as three separate runs: "(", "분류모델", ")".
# Each InsertAfter call on the paragraph's TextRange appends a brand-new
# run just before the paragraph's endParaRPr, inheriting the run
# formatting (bold, 24pt, white fill) of the preceding run - matching
# the look of the existing "예측" run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$titleRange = $shape.TextFrame.TextRange

$openParen = $titleRange.InsertAfter("(")
$classModel = $titleRange.InsertAfter("분류모델")
$closeParen = $titleRange.InsertAfter(")")
